$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DocEntry list (header stays "DocEntry"; rows 2-21 replace the old
# numeric SAP doc-entry values with the new order/request codes). The old
# rows carried a highlight style (yellow fill + bottom border on some
# cells) that is no longer needed for this data set, so clear cell
# formatting first.
$ws.Range("A1:A17").ClearFormats()

$values = @(
    "DocEntry",
    "2025-9539-23APR-3",
    "2025-121390",
    "2025-121388",
    "2025-121735",
    "2025-121735-A1",
    "2025-122352",
    "2025-122352-A2",
    "2025-10347-20JUN-2",
    "2025-10347-20JUN-3",
    "2025-123533",
    "REQ-13",
    "2025-10582-19JUL-1",
    "2025-124983",
    "2025-128558",
    "2025-128547",
    "2025-128536",
    "2025-128584",
    "2025-11704-20OCT-1",
    "2025-129624",
    "2025-129986-A1"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Size column A to fit the new (longer) text values.
$ws.Columns("A").ClearFormats()
$ws.Columns("A").AutoFit()

# Match the author's final selection/cursor position.
$ws.Range("C10").Select()
